# Automatische test-sync: 2025-06-19 21:15:50
# Append 4 new inbound mail rows to "Logs" and their matching category
# roll-up rows to "Dashboard", then widen the chart series ranges and
# the conditional-formatting ranges so they cover the new rows too.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# ---------------------------------------------------------------------
# Logs sheet — rows 3..6
# ---------------------------------------------------------------------
$logs.Range("A3").Value = "Wat zijn jullie openingstijden?"
$logs.Range("B3").Value = "mailmind.test@zohomail.eu"
$logs.Range("C3").Value = "Hallo, ik zou graag willen weten wat jullie openingstijden zijn. Dank je wel!"
$logs.Range("D3").Value = "Openingstijden / Locatie"
$logs.Range("E3").Value = "Beste,`nBedankt voor je interesse. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 18:00 uur en op zaterdag van 10:00 tot 16:00 uur. Op zondag zijn wij gesloten. Mocht je verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf]"
$logs.Range("F3").Value = "2025-06-19 21:12:26"
$logs.Range("G3").Value = "Ja"

$logs.Range("A4").Value = "Factuur verzoek"
$logs.Range("B4").Value = "mailmind.test@zohomail.eu"
$logs.Range("C4").Value = "Kunt u mij de factuur van mijn laatste bestelling toesturen?"
$logs.Range("D4").Value = "Factuur / Administratie"
$logs.Range("F4").Value = "2025-06-19 21:13:10"
$logs.Range("G4").Value = "Nee"

$logs.Range("A5").Value = "Is product X op voorraad?"
$logs.Range("B5").Value = "mailmind.test@zohomail.eu"
$logs.Range("C5").Value = "Ik ben geïnteresseerd in product X. Is dit momenteel op voorraad?"
$logs.Range("D5").Value = "Productinformatie"
$logs.Range("F5").Value = "2025-06-19 21:14:09"
$logs.Range("G5").Value = "Nee"

$logs.Range("A6").Value = "Probleem met inloggen"
$logs.Range("B6").Value = "mailmind.test@zohomail.eu"
$logs.Range("C6").Value = "Ik kan niet inloggen op mijn account. Kunnen jullie dit oplossen?"
$logs.Range("D6").Value = "IT / Technisch probleem"
$logs.Range("F6").Value = "2025-06-19 21:15:14"
$logs.Range("G6").Value = "Nee"

# Conditional formatting on Logs must now cover the new rows too.
$logs.Range("D2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D6"))
$logs.Range("G2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G6"))

# ---------------------------------------------------------------------
# Dashboard sheet — rows 3..6 (one per new category)
# ---------------------------------------------------------------------
$dash.Range("A3").Value = "Openingstijden / Locatie"
$dash.Range("B3").Value = 1

$dash.Range("A4").Value = "Factuur / Administratie"
$dash.Range("B4").Value = 1

$dash.Range("A5").Value = "Productinformatie"
$dash.Range("B5").Value = 1

$dash.Range("A6").Value = "IT / Technisch probleem"
$dash.Range("B6").Value = 1

# ---------------------------------------------------------------------
# Chart1 on the Dashboard sheet — widen category/value series ranges
# from the single row 2 to the full 2..6 block.
# ---------------------------------------------------------------------
$chart = $dash.ChartObjects(1).Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$6"
$series.Values = "='Dashboard'!`$B`$2:`$B`$6"

Write-Output "sync applied"
